# Update cryptos list with latest scraped values (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.275.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.54%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.28%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.23%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4701"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.27%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3930"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.87%  "

# Row 9 - OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07985"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.92%  "

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "

# Row 12 - Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.13%  "

# Row 13 - WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.872.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.989"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.278"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.34%  "

# Row 16 - BinanceUSD (was Litecoin)
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "

# Row 17 - Litecoin (was BinanceUSD)
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.52%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001044"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19 - TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06584"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.89%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.55%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.29%  "

# Row 22 - WrappedBTC
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.289.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.60%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.449"
$ws.Range("D23").Style = "Normal"

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.05%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.297"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.107.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.36%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.70%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.160"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.75%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.66%  "

# Row 31 - BitcoinCash
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.09%  "

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9796"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.29%  "

# Row 33 - Stellar
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09484"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.81%  "

# Row 34 - HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.583"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.380"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.90%  "

# Row 36 - Filecoin
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.357"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.36%  "

# Row 37 - VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02273"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.14%  "

# Row 38 - Hedera
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06084"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.04%  "

# Row 39 - FraxShare
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.459"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.88%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.180"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.31%  "

# Row 41 - TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5962"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.21%  "

# Row 42 - Frax
$ws.Range("E42").Value = "  +0.13%  "

# Row 43 - Algorand
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1876"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.77%  "

# Row 45 - WEMIXTOKEN
$ws.Range("E45").Value = "  +4.95%  "

# Row 46 - Decentraland
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5614"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.66%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.33%  "

# Row 48 - NEARProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.967"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.41%  "

# Row 49 - Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06895"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.19%  "

# Row 50 - Quant
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.07%  "

# Row 51 - RenderToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.017"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.02%  "
